$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column cells remain text (they are stored as text in the source data,
# e.g. "28.880.07" using dots as thousands separators -- not valid numbers).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.880.07"
$ws.Range("E2").Value = "  -0.31%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.877.61"
$ws.Range("E3").Value = "  -0.94%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.84"
$ws.Range("E5").Value = "  -0.75%  "
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("E7").Value = "  -0.36%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3864"
$ws.Range("E8").Value = "  -0.74%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07848"
$ws.Range("E9").Value = "  -0.91%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9851"
$ws.Range("E10").Value = "  -2.55%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.75"
$ws.Range("E11").Value = "  -0.89%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.899.92"
$ws.Range("E12").Value = "  -0.75%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.981"
$ws.Range("E13").Value = "  -1.77%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.645"
$ws.Range("E14").Value = "  -1.87%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06975"
$ws.Range("E15").Value = "  -0.03%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.05"
$ws.Range("E16").Value = "  -0.21%  "
$ws.Range("E17").Value = "  +0.06%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009950"
$ws.Range("E18").Value = "  -1.33%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.92"
$ws.Range("E19").Value = "  -2.06%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.002"
$ws.Range("E20").Value = "  -0.48%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "28.893.26"
$ws.Range("E21").Value = "  -0.20%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.241"
$ws.Range("E22").Value = "  -2.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.93"
$ws.Range("E23").Value = "  -1.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.103"
$ws.Range("E24").Value = "  +2.00%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "156.24"
$ws.Range("E25").Value = "  +0.65%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "19.33"
$ws.Range("E26").Value = "  -1.94%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "5.981"
$ws.Range("E27").Value = "  +1.68%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "117.39"
$ws.Range("E28").Value = "  -1.44%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.914"
$ws.Range("E29").Value = "  -3.79%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09361"
$ws.Range("E30").Value = "  -0.16%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.9006"
$ws.Range("E31").Value = "  -3.64%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.258"
$ws.Range("E32").Value = "  -1.82%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.315"
$ws.Range("E33").Value = "  -2.34%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.253"
$ws.Range("E34").Value = "  -0.04%  "
$ws.Range("E35").Value = "  +1.03%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05741"
$ws.Range("E36").Value = "  -1.23%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02073"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.002"
$ws.Range("E38").Value = "  -0.22%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.629"
$ws.Range("E39").Value = "  -5.88%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5638"
$ws.Range("E40").Value = "  -2.52%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1763"
$ws.Range("E41").Value = "  -2.61%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.690"
$ws.Range("E42").Value = "  -2.62%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.274"
$ws.Range("E43").Value = "  +3.59%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "11.94"
$ws.Range("E44").Value = "  +0.48%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5339"
$ws.Range("E45").Value = "  -1.84%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.07042"
$ws.Range("E46").Value = "  -1.90%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.841"
$ws.Range("E47").Value = "  -1.16%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.546"
$ws.Range("E48").Value = "  +2.16%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "112.78"
$ws.Range("E49").Value = "  -0.52%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.069"
$ws.Range("E50").Value = "  -5.13%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "70.72"
$ws.Range("E51").Value = "  -0.83%  "
